# Update countries & provincias Spain
# - Re-rank a handful of countries whose case counts crossed (name swaps in col A)
# - Refresh the "last updated" footer timestamp
# - Refresh the daily COVID-19 statistics (cols B-H) for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country ranking swaps (column A) ---
$ws.Range("A47").Value = "Japon"
$ws.Range("A48").Value = "Polonia"
$ws.Range("A94").Value = "Guayana Francesa"
$ws.Range("A95").Value = "Grecia"
$ws.Range("A105").Value = "Zimbabue"
$ws.Range("A106").Value = "Namibia"
$ws.Range("A154").Value = "Togo"
$ws.Range("A155").Value = "Liberia"
$ws.Range("A158").Value = "Trinidad yTobago"
$ws.Range("A159").Value = "Principado de Andorra"
$ws.Range("A214").Value = "Montserrat"
$ws.Range("A215").Value = "Islas Malvinas"

# --- Updated "last refreshed" footer (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Agosto de 2020 a las 01:01"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5913180
$ws.Range("C4").Value = 39034
$ws.Range("D4").Value = 3207475
$ws.Range("E4").Value = 2524646
$ws.Range("G4").Value = 455
$ws.Range("H4").Value = 181059
$ws.Range("D5").Value = 2778709
$ws.Range("E5").Value = 728843
$ws.Range("B9").Value = 600438
$ws.Range("C9").Value = 6112
$ws.Range("D9").Value = 407301
$ws.Range("E9").Value = 165324
$ws.Range("G9").Value = 150
$ws.Range("H9").Value = 27813
$ws.Range("B15").Value = 350867
$ws.Range("C15").Value = 8713
$ws.Range("E15").Value = 86712
$ws.Range("G15").Value = 381
$ws.Range("H15").Value = 7366
$ws.Range("B34").Value = 97478
$ws.Range("C34").Value = 138
$ws.Range("D34").Value = 66817
$ws.Range("E34").Value = 25381
$ws.Range("G34").Value = 18
$ws.Range("H34").Value = 5280
$ws.Range("B35").Value = 91608
$ws.Range("C35").Value = 447
$ws.Range("D35").Value = 61558
$ws.Range("E35").Value = 28477
$ws.Range("G35").Value = 6
$ws.Range("H35").Value = 1573
$ws.Range("B36").Value = 87485
$ws.Range("C36").Value = 585
$ws.Range("D36").Value = 62185
$ws.Range("E36").Value = 23394
$ws.Range("G36").Value = 14
$ws.Range("H36").Value = 1906
$ws.Range("B44").Value = 68533
$ws.Range("C44").Value = 345
$ws.Range("D44").Value = 57735
$ws.Range("E44").Value = 8187
$ws.Range("G44").Value = 17
$ws.Range("H44").Value = 2611
$ws.Range("B47").Value = 62507
$ws.Range("C47").Value = 760
$ws.Range("D47").Value = 49340
$ws.Range("E47").Value = 11986
$ws.Range("H47").Value = 1181
$ws.Range("B48").Value = 62310
$ws.Range("C48").Value = 548
$ws.Range("D48").Value = 42448
$ws.Range("E48").Value = 17902
$ws.Range("G48").Value = 5
$ws.Range("H48").Value = 1960
$ws.Range("B87").Value = 12903
$ws.Range("C87").Value = 67
$ws.Range("D87").Value = 6538
$ws.Range("E87").Value = 5547
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 818
$ws.Range("B88").Value = 11148
$ws.Range("C88").Value = 66
$ws.Range("D88").Value = 10208
$ws.Range("E88").Value = 660
$ws.Range("B93").Value = 9076
$ws.Range("C93").Value = 109
$ws.Range("D93").Value = 7928
$ws.Range("E93").Value = 1093
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 55
$ws.Range("B94").Value = 8875
$ws.Range("C94").Value = 78
$ws.Range("D94").Value = 8363
$ws.Range("E94").Value = 456
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 56
$ws.Range("B95").Value = 8819
$ws.Range("C95").Value = 155
$ws.Range("D95").Value = 3804
$ws.Range("E95").Value = 4773
$ws.Range("H95").Value = 242
$ws.Range("B97").Value = 8409
$ws.Range("C97").Value = 21
$ws.Range("D97").Value = 6959
$ws.Range("E97").Value = 1397
$ws.Range("B103").Value = 6928
$ws.Range("C103").Value = 23
$ws.Range("D103").Value = 6282
$ws.Range("E103").Value = 488
$ws.Range("B104").Value = 6912
$ws.Range("C104").Value = 133
$ws.Range("D104").Value = 4297
$ws.Range("E104").Value = 2588
$ws.Range("B105").Value = 6070
$ws.Range("C105").Value = 140
$ws.Range("D105").Value = 4950
$ws.Range("E105").Value = 965
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 155
$ws.Range("B106").Value = 6030
$ws.Range("C106").Value = 176
$ws.Range("D106").Value = 2563
$ws.Range("E106").Value = 3411
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 56
$ws.Range("B107").Value = 5419
$ws.Range("C107").Value = 5
$ws.Range("D107").Value = 3059
$ws.Range("E107").Value = 2191
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 169
$ws.Range("B115").Value = 4304
$ws.Range("C115").Value = 79
$ws.Range("D115").Value = 2936
$ws.Range("E115").Value = 1283
$ws.Range("B116").Value = 3979
$ws.Range("C116").Value = 129
$ws.Range("D116").Value = 1742
$ws.Range("E116").Value = 2159
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 78
$ws.Range("B119").Value = 3532
$ws.Range("C119").Value = 23
$ws.Range("D119").Value = 2599
$ws.Range("E119").Value = 896
$ws.Range("B136").Value = 2222
$ws.Range("C136").Value = 51
$ws.Range("D136").Value = 877
$ws.Range("E136").Value = 1245
$ws.Range("G136").Value = 4
$ws.Range("H136").Value = 100
$ws.Range("B140").Value = 1997
$ws.Range("C140").Value = 5
$ws.Range("D140").Value = 1557
$ws.Range("E140").Value = 371
$ws.Range("B154").Value = 1295
$ws.Range("C154").Value = 18
$ws.Range("D154").Value = 914
$ws.Range("E154").Value = 354
$ws.Range("H154").Value = 27
$ws.Range("B155").Value = 1290
$ws.Range("C155").Value = 4
$ws.Range("D155").Value = 819
$ws.Range("E155").Value = 389
$ws.Range("H155").Value = 82
$ws.Range("D157").Value = 1084
$ws.Range("E157").Value = 19
$ws.Range("B158").Value = 1099
$ws.Range("C158").Value = 92
$ws.Range("D158").Value = 165
$ws.Range("E158").Value = 919
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 15
$ws.Range("B159").Value = 1060
$ws.Range("C159").Value = 15
$ws.Range("D159").Value = 877
$ws.Range("E159").Value = 130
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 53
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
